$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert two new blank rows right below the header (before the existing data)
$ws.Rows.Item(2).Resize(2,1).Insert() | Out-Null

# Grow the table to include the two new rows
$lo.Resize($ws.Range("A1:C8")) | Out-Null

# Fill the new rows with the new drivers (CARLOS first so it lands at
# shared-string index 8, then RENATO at index 9 - matches canonical order)
$ws.Range("A3").Value = "CARLOS ALBERTO DE MENEZES"
$ws.Range("B3").Value = 666
$ws.Range("C3").Value = 60

$ws.Range("A2").Value = "RENATO AQUINO DE PIN"
$ws.Range("B2").Value = 777
$ws.Range("C2").Value = 70

# Move the special underline style from the old A2 (now A4) to the new A3
$ws.Range("A4").Font.Underline = $false
$ws.Range("A3").Font.Underline = $true

# Update selection to A3
$ws.Range("A3").Select() | Out-Null
